$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Primera" and "Segunda" quality corazón de apio entries
# were swapped between row 3 and row 5 (row 3 now holds the "Segunda" data
# that used to be in row 5, and vice versa).

# New values for row 3 (previously held by row 5)
$ws.Range("D3").Value = 44377
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 550
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = 2364
$ws.Range("P3").Value = 394

# New values for row 5 (previously held by row 3)
$ws.Range("D5").Value = 44267
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1800
$ws.Range("M5").Value = 1650
$ws.Range("P5").Value = 275
